$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 96 - 2025-09-17, 四方坪站
$ws.Range("A96").Value = 45917
$ws.Range("B96").Value = "四方坪站"
$ws.Range("C96").Value = 11856.32
$ws.Range("D96").Value = 9661.43
$ws.Range("E96").Value = 4143.76
$ws.Range("F96").Value = 484

# Row 97 - 2025-09-17, 高岭站
$ws.Range("A97").Value = 45917
$ws.Range("B97").Value = "高岭站"
$ws.Range("C97").Value = 5535.89
$ws.Range("D97").Value = 4324.51
$ws.Range("E97").Value = 1388.51
$ws.Range("F97").Value = 185

# Update the view: scroll to the new rows and select I93 (matches author's cursor position)
$ws.Range("I93").Select() | Out-Null
